# CheckList.xlsx update: the two team members previously listed as
# "Vo Gia Huy" (Assigned to / Update columns) are replaced by two
# different people - "Hoang Viet Anh" (Assigned to, column E) and
# "Dao Vinh Phat" (Update, column G) - across every task row (3-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$assignedTo = "Hoang Viet Anh"
$updatedBy  = "Dao Vinh Phat"

for ($row = 3; $row -le 8; $row++) {
    $ws.Cells.Item($row, 5).Value = $assignedTo   # column E - Assigned to
    $ws.Cells.Item($row, 7).Value = $updatedBy    # column G - Update
}

# Widen the Assigned-to / Date-modify / Update columns to fit the new names
# (values chosen so the saved column width lands as close as possible to
# the author's original 20.77734375 / 19.109375 / 17.6640625).
$ws.Columns.Item(5).ColumnWidth = 120 / 6
$ws.Columns.Item(6).ColumnWidth = 110 / 6
$ws.Columns.Item(7).ColumnWidth = 101 / 6

# Move the active selection (reflects where the author last clicked).
$ws.Range("F11").Select()

Write-Output "CheckList updated"
